# TAKEN 1 & 2
# Rename the worksheet/tab to match the newly downloaded response export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "TASK0-Your responses.downloadlo"

# The "Marking Status" column (K) moved from "Marked" to "Posted" for every
# answer row.
$ws.Range("K2").Value = "Posted"
$ws.Range("K3").Value = "Posted"
$ws.Range("K4").Value = "Posted"

# The "Answer" column (G) held mis-typed / mis-scaled values (an integer
# "1357" formatted with a thousands separator, and two answers stored as
# text). Replace them with the real numeric answers.
$ws.Range("G2").ClearFormats()
$ws.Range("G2").Value = 1.357
$ws.Range("G3").Value = 0.861
$ws.Range("G4").Value = 0.587

$wb.Save()
